$d = $word.ActiveDocument

# 1. Merge the split "T" / "his is an R Markdown document..." runs into a single run.
$d.Content.Find.Execute("This is an R Markdown document. Markdown is a simple formatting syntax for authoring HTML, PDF, and MS Word documents. For more details on using R Markdown see ", $true, $false, $false, $false, $false, $true, 1, $false, "This is an R Markdown document. Markdown is a simple formatting syntax for authoring HTML, PDF, and MS Word documents. For more details on using R Markdown see ", 2) | Out-Null

# 2. Merge the split "s" / "ummary" runs (FunctionTok style) into a single run.
$d.Content.Find.Execute("summary", $true, $false, $false, $false, $false, $true, 1, $false, "summary", 2) | Out-Null

# 3. Merge the split VerbatimChar runs "##  Min.   : 4.0   Min." / "   :  2.00  ".
$d.Content.Find.Execute("##  Min.   : 4.0   Min.   :  2.00  ", $true, $false, $false, $false, $false, $true, 1, $false, "##  Min.   : 4.0   Min.   :  2.00  ", 2) | Out-Null

# 4. Merge the split "Note tha" / "t the " runs into a single run.
$d.Content.Find.Execute("Note that the ", $true, $false, $false, $false, $false, $true, 1, $false, "Note that the ", 2) | Out-Null

# 5. "Image Caption" style: no longer italic, and a smaller (9pt) font size --
#    this is what lets captions sit comfortably under figures.
$imageCaption = $d.Styles.Item("ImageCaption")
$imageCaption.Font.Italic = 0
$imageCaption.Font.Size = 9

# 6. "Figure" style: center the figure paragraph on the page.
$figure = $d.Styles.Item("Figure")
$figure.ParagraphFormat.Alignment = 1

